# "Redesign of the HTM with movement based on gravity model."
# On the "Setline Survey" sheet, the years 1987-1992 (rows 12:17) are removed
# entirely (a real row delete, so rows below shift up), and the numeric data
# columns (B:I) for all remaining data rows get a "0.00" number format.
# The active selection ends up on F37 (left over from the editing session,
# even though it now falls outside the shrunk used range).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setline Survey")

# Remove the rows for 1987-1992 (original rows 12 through 17); this shifts
# every row below up by six, turning old row 18 (1993) into new row 12, and
# old row 37 (2012) into new row 31.
$ws.Rows("12:17").Delete()

# Apply a two-decimal numeric format to the data columns (B:I) for every
# data row that remains (rows 2 through 31).
$ws.Range("B2:I31").NumberFormat = "0.00"

# Leave the selection where the editing session left it.
$ws.Range("F37").Select()
